$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$newRow = $t.Rows.Add()
$newRow.Cells.Item(1).Range.Text = "Merci !"
$newRow.Cells.Item(2).Range.Text = "Merci aussi à vous. Nous à votre disposition pour toutes préoccupations"
